# Insert a new column before C. This shifts the existing C,D,E,F columns
# (and all their cell content / formulas) one place to the right, becoming
# D,E,F,G respectively - matching the target diff exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C:C").Insert()

# New column C should look like column B (same header style on row 1,
# same body style on rows 2:11). Copy the formatting over, then fill in
# the new "IconSource" header + per-row icon class values.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("B2:B11").Copy() | Out-Null
$ws.Range("C2:C11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New column's width (as close as this host's ColumnWidth quantization
# allows to the authored 16.42578125 character width).
$ws.Columns("C:C").ColumnWidth = 15.67

$ws.Range("C1").Value = "IconSource"

$icons = @(
  "fas fa-home",
  "fas fa-shopping-basket",
  "fas fa-users",
  "fas fa-credit-card",
  "fas fa-chart-line",
  "fas fa-hourglass-start",
  "fas fa-cart-arrow-down",
  "fas fa-warehouse",
  "fas fa-registered",
  "fas fa-user"
)

for ($i = 0; $i -lt $icons.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 3).Value = $icons[$i]
}

# The SQL-building formula (now in column G after the insert) needs the
# new IconSource column appended as an extra CONCATENATE argument. Rewrite
# row 2 directly and let the rest of the column pick it up as one shared
# formula, same as the original authoring pattern.
$ws.Range("G2").Formula = '=CONCATENATE("PERFORM ""SchSysConfig"".""Func_TblAppObject_MenuGroup_SET""(varSystemLoginSession, null, null, null, varInstitutionBranchID, varBaseCurrencyID, ''", B2, "''::varchar, ''", C2, "''::varchar);")'
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G3:G11").PasteSpecial(-4123) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("G15").Select()
